$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a text value to a cell while preventing Excel from
# auto-converting numeric-looking strings into real numbers, and without
# leaving a lasting style/number-format change behind on the cell.
function Set-TextValue($cell, $text) {
    $originalStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $originalStyle
}

$ws.Range("D2").Value = '88.279.65'
$ws.Range("E2").Value = '  -1.96%  '

$ws.Range("D3").Value = '3.085.41'
$ws.Range("E3").Value = '  -3.90%  '

$ws.Range("E4").Value = '  -0.06%  '

Set-TextValue $ws.Range("D5") '209.60'
$ws.Range("E5").Value = '  -4.13%  '

Set-TextValue $ws.Range("D6") '623.36'
$ws.Range("E6").Value = '  -0.44%  '

Set-TextValue $ws.Range("D7") '0.372'
$ws.Range("E7").Value = '  -4.26%  '

Set-TextValue $ws.Range("D8") '0.820'
$ws.Range("E8").Value = '  +17.65%  '

$ws.Range("E9").Value = '  +0.04%  '

$ws.Range("D10").Value = '3.083.80'
$ws.Range("E10").Value = '  -3.78%  '

$ws.Range("E11").Value = '  +3.67%  '

$ws.Range("E12").Value = '  +0.05%  '

$ws.Range("E13").Value = '  -6.57%  '

Set-TextValue $ws.Range("D14") '5.28'
$ws.Range("E14").Value = '  -2.21%  '

$ws.Range("D15").Value = '87.978.46'
$ws.Range("E15").Value = '  -2.09%  '

$ws.Range("D16").Value = '3.653.65'
$ws.Range("E16").Value = '  -3.79%  '

Set-TextValue $ws.Range("D17") '31.62'
$ws.Range("E17").Value = '  -5.05%  '

$ws.Range("D18").Value = '3.077.33'
$ws.Range("E18").Value = '  -4.15%  '

$ws.Range("E19").Value = '  -5.61%  '

$ws.Range("E20").Value = '  -10.84%  '

Set-TextValue $ws.Range("D21") '13.09'
$ws.Range("E21").Value = '  -2.67%  '

Set-TextValue $ws.Range("D22") '421.07'
$ws.Range("E22").Value = '  -4.72%  '

Set-TextValue $ws.Range("D23") '8.17'
$ws.Range("E23").Value = '  -4.76%  '

$ws.Range("E24").Value = '  -4.41%  '

$ws.Range("E25").Value = '  +6.32%  '

Set-TextValue $ws.Range("D26") '81.64'
$ws.Range("E26").Value = '  -0.24%  '

Set-TextValue $ws.Range("D27") '11.49'
$ws.Range("E27").Value = '  -1.53%  '

$ws.Range("D28").Value = '3.244.18'
$ws.Range("E28").Value = '  -3.36%  '

$ws.Range("E29").Value = '  -0.01%  '

$ws.Range("E30").Value = '  +8.63%  '

$ws.Range("E31").Value = '  -1.55%  '

Set-TextValue $ws.Range("D32") '8.03'
$ws.Range("E32").Value = '  -5.35%  '

Set-TextValue $ws.Range("D33") '504.06'
$ws.Range("E33").Value = '  -6.48%  '

Set-TextValue $ws.Range("D34") '3.54'
$ws.Range("E34").Value = '  -13.09%  '

Set-TextValue $ws.Range("D35") '6.60'
$ws.Range("E35").Value = '  -4.41%  '

$ws.Range("E36").Value = '  -4.81%  '

$ws.Range("E37").Value = '  -4.47%  '

Set-TextValue $ws.Range("D38") '22.20'
$ws.Range("E38").Value = '  -1.58%  '

Set-TextValue $ws.Range("D39") '0.132'
$ws.Range("E39").Value = '  +3.51%  '

$ws.Range("E40").Value = '  -0.45%  '

$ws.Range("E41").Value = '  +0.25%  '

Set-TextValue $ws.Range("D43") '149.35'
$ws.Range("E43").Value = '  -0.28%  '

$ws.Range("B44").Value = 'PolygonEcosystemToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextValue $ws.Range("D44") '0.358'
$ws.Range("E44").Value = '  -3.75%  '

$ws.Range("B45").Value = 'Stellar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range("D45") '0.135'
$ws.Range("E45").Value = '  +9.08%  '

Set-TextValue $ws.Range("D46") '1.81'
$ws.Range("E46").Value = '  -5.62%  '

Set-TextValue $ws.Range("D47") '43.47'
$ws.Range("E47").Value = '  +0.36%  '

Set-TextValue $ws.Range("D48") '0.0660'
$ws.Range("E48").Value = '  +10.11%  '

Set-TextValue $ws.Range("D49") '0.702'
$ws.Range("E49").Value = '  -4.88%  '

Set-TextValue $ws.Range("D50") '155.94'
$ws.Range("E50").Value = '  -9.34%  '

$ws.Range("E51").Value = '  -5.57%  '
